$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 11 for the Xbee USB programmer entry; everything
# below (rows 11-41) shifts down by one (new rows 12-42).
$ws.Rows.Item(11).Insert()

$ws.Range("B11").Value = "WRL-11812"
$ws.Range("C11").Value = "SparkFun XBee Explorer USB (used to program the Xbee module)`nhttp://www.sparkfun.com/products/11812"
$ws.Range("D11").Value = "Sparkfun"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 24.95
$ws.Range("G11").Formula = "=F11*E11"
$ws.Rows.Item(11).RowHeight = 30

# Restore the originally-selected cell as recorded in the saved workbook.
$ws.Range("C12").Select()
